$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 4.713280333333334
$ws.Range("H2").Value = 14.139841
$ws.Range("I2").Value = 0.6942627767023021
$ws.Range("J2").Value = 0.6942627767023022
$ws.Range("M2").Value = 1.711937666666667
$ws.Range("N2").Value = 5.135813
$ws.Range("O2").Value = 0.09827532014408574
$ws.Range("P2").Value = 0.09827532014408574
$ws.Range("Q2").Value = 8.068842136192556
$ws.Range("R2").Value = 72.61957922573301
$ws.Range("S2").Value = 0.06822889664454065
$ws.Range("T2").Value = 0.06822889664454065
$ws.Range("G3").Value = 4.713280333333334
$ws.Range("H3").Value = 14.139841
$ws.Range("I3").Value = 0.6942627767023021
$ws.Range("J3").Value = 0.6942627767023022
$ws.Range("O3").Value = 0.3329367223581701
$ws.Range("P3").Value = 0.3329367223581701
$ws.Range("Q3").Value = 27.33558995392511
$ws.Range("R3").Value = 246.020309585326
$ws.Range("S3").Value = 0.2311455733305466
$ws.Range("T3").Value = 0.2311455733305466
$ws.Range("G4").Value = 4.713280333333334
$ws.Range("H4").Value = 14.139841
$ws.Range("I4").Value = 0.6942627767023021
$ws.Range("J4").Value = 0.6942627767023022
$ws.Range("M4").Value = 3.605537
$ws.Range("N4").Value = 10.816611
$ws.Range("O4").Value = 0.2069790915087912
$ws.Range("P4").Value = 0.2069790915087912
$ws.Range("Q4").Value = 16.99390663320566
$ws.Range("R4").Value = 152.945159698851
$ws.Range("S4").Value = 0.1436978787902132
$ws.Range("T4").Value = 0.1436978787902133
$ws.Range("G5").Value = 4.713280333333334
$ws.Range("H5").Value = 14.139841
$ws.Range("I5").Value = 0.6942627767023021
$ws.Range("J5").Value = 0.6942627767023022
$ws.Range("M5").Value = 6.302642666666666
$ws.Range("N5").Value = 18.907928
$ws.Range("O5").Value = 0.3618088659889531
$ws.Range("P5").Value = 0.361808865988953
$ws.Range("Q5").Value = 29.70612172882755
$ws.Range("R5").Value = 267.355095559448
$ws.Range("S5").Value = 0.2511904279370016
$ws.Range("T5").Value = 0.2511904279370016
$ws.Range("I6").Value = 0.1843120478188439
$ws.Range("J6").Value = 0.1843120478188439
$ws.Range("M6").Value = 1.711937666666667
$ws.Range("N6").Value = 5.135813
$ws.Range("O6").Value = 0.09827532014408574
$ws.Range("P6").Value = 0.09827532014408574
$ws.Range("Q6").Value = 2.142106515796
$ws.Range("R6").Value = 19.278958642164
$ws.Range("S6").Value = 0.01811332550580892
$ws.Range("T6").Value = 0.01811332550580893
$ws.Range("I7").Value = 0.1843120478188439
$ws.Range("J7").Value = 0.1843120478188439
$ws.Range("O7").Value = 0.3329367223581701
$ws.Range("P7").Value = 0.3329367223581701
$ws.Range("S7").Value = 0.06136424909192819
$ws.Range("T7").Value = 0.0613642490919282
$ws.Range("I8").Value = 0.1843120478188439
$ws.Range("J8").Value = 0.1843120478188439
$ws.Range("M8").Value = 3.605537
$ws.Range("N8").Value = 10.816611
$ws.Range("O8").Value = 0.2069790915087912
$ws.Range("P8").Value = 0.2069790915087912
$ws.Range("Q8").Value = 4.511521915212
$ws.Range("R8").Value = 40.603697236908
$ws.Range("S8").Value = 0.03814874021166918
$ws.Range("T8").Value = 0.03814874021166919
$ws.Range("I9").Value = 0.1843120478188439
$ws.Range("J9").Value = 0.1843120478188439
$ws.Range("M9").Value = 6.302642666666666
$ws.Range("N9").Value = 18.907928
$ws.Range("O9").Value = 0.3618088659889531
$ws.Range("P9").Value = 0.361808865988953
$ws.Range("Q9").Value = 7.886345505376
$ws.Range("R9").Value = 70.97710954838401
$ws.Range("S9").Value = 0.0666857330094376
$ws.Range("T9").Value = 0.0666857330094376
$ws.Range("G10").Value = 0.6263116666666667
$ws.Range("H10").Value = 1.878935
$ws.Range("I10").Value = 0.09225525452111802
$ws.Range("J10").Value = 0.09225525452111803
$ws.Range("M10").Value = 1.711937666666667
$ws.Range("N10").Value = 5.135813
$ws.Range("O10").Value = 0.09827532014408574
$ws.Range("P10").Value = 0.09827532014408574
$ws.Range("Q10").Value = 1.072206533239444
$ws.Range("R10").Value = 9.649858799155
$ws.Range("S10").Value = 0.009066414673036987
$ws.Range("T10").Value = 0.009066414673036989
$ws.Range("G11").Value = 0.6263116666666667
$ws.Range("H11").Value = 1.878935
$ws.Range("I11").Value = 0.09225525452111802
$ws.Range("J11").Value = 0.09225525452111803
$ws.Range("O11").Value = 0.3329367223581701
$ws.Range("P11").Value = 0.3329367223581701
$ws.Range("Q11").Value = 3.632416850378889
$ws.Range("R11").Value = 32.69175165341
$ws.Range("S11").Value = 0.03071516206057978
$ws.Range("T11").Value = 0.03071516206057979
$ws.Range("G12").Value = 0.6263116666666667
$ws.Range("H12").Value = 1.878935
$ws.Range("I12").Value = 0.09225525452111802
$ws.Range("J12").Value = 0.09225525452111803
$ws.Range("M12").Value = 3.605537
$ws.Range("N12").Value = 10.816611
$ws.Range("O12").Value = 0.2069790915087912
$ws.Range("P12").Value = 0.2069790915087912
$ws.Range("Q12").Value = 2.258189887698333
$ws.Range("R12").Value = 20.323708989285
$ws.Range("S12").Value = 0.01909490876769331
$ws.Range("T12").Value = 0.01909490876769331
$ws.Range("G13").Value = 0.6263116666666667
$ws.Range("H13").Value = 1.878935
$ws.Range("I13").Value = 0.09225525452111802
$ws.Range("J13").Value = 0.09225525452111803
$ws.Range("M13").Value = 6.302642666666666
$ws.Range("N13").Value = 18.907928
$ws.Range("O13").Value = 0.3618088659889531
$ws.Range("P13").Value = 0.361808865988953
$ws.Range("Q13").Value = 3.947418632964445
$ws.Range("R13").Value = 35.52676769668
$ws.Range("S13").Value = 0.03337876901980794
$ws.Range("T13").Value = 0.03337876901980794
$ws.Range("G14").Value = 0.1980316666666667
$ws.Range("H14").Value = 0.594095
$ws.Range("I14").Value = 0.02916992095773596
$ws.Range("J14").Value = 0.02916992095773596
$ws.Range("M14").Value = 1.711937666666667
$ws.Range("N14").Value = 5.135813
$ws.Range("O14").Value = 0.09827532014408574
$ws.Range("P14").Value = 0.09827532014408574
$ws.Range("Q14").Value = 0.3390178693594444
$ws.Range("R14").Value = 3.051160824235
$ws.Range("S14").Value = 0.002866683320699177
$ws.Range("T14").Value = 0.002866683320699178
$ws.Range("G15").Value = 0.1980316666666667
$ws.Range("H15").Value = 0.594095
$ws.Range("I15").Value = 0.02916992095773596
$ws.Range("J15").Value = 0.02916992095773596
$ws.Range("O15").Value = 0.3329367223581701
$ws.Range("P15").Value = 0.3329367223581701
$ws.Range("Q15").Value = 1.148523333018889
$ws.Range("R15").Value = 10.33670999717
$ws.Range("S15").Value = 0.009711737875115502
$ws.Range("T15").Value = 0.009711737875115504
$ws.Range("G16").Value = 0.1980316666666667
$ws.Range("H16").Value = 0.594095
$ws.Range("I16").Value = 0.02916992095773596
$ws.Range("J16").Value = 0.02916992095773596
$ws.Range("M16").Value = 3.605537
$ws.Range("N16").Value = 10.816611
$ws.Range("O16").Value = 0.2069790915087912
$ws.Range("P16").Value = 0.2069790915087912
$ws.Range("Q16").Value = 0.7140105013383333
$ws.Range("R16").Value = 6.426094512044999
$ws.Range("S16").Value = 0.006037563739215436
$ws.Range("T16").Value = 0.006037563739215437
$ws.Range("G17").Value = 0.1980316666666667
$ws.Range("H17").Value = 0.594095
$ws.Range("I17").Value = 0.02916992095773596
$ws.Range("J17").Value = 0.02916992095773596
$ws.Range("M17").Value = 6.302642666666666
$ws.Range("N17").Value = 18.907928
$ws.Range("O17").Value = 0.3618088659889531
$ws.Range("P17").Value = 0.361808865988953
$ws.Range("Q17").Value = 1.248122831684445
$ws.Range("R17").Value = 11.23310548516
$ws.Range("S17").Value = 0.01055393602270584
$ws.Range("T17").Value = 0.01055393602270584
